# Costas 2nd order loop filter
# Applies the changes described by the diff to 300_iir_llpf_calc.xlsx

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "LoopFilter LPF" sheet: change Gain (B16) from 8 to 1
#    (dependent formulas B20/B21 recalc automatically)
# ---------------------------------------------------------------------
$loopFilter = $wb.Worksheets.Item("LoopFilter LPF")
$loopFilter.Range("B16").Value = 1

# ---------------------------------------------------------------------
# 2. Insert a new blank worksheet named "Sheet1" right after
#    "LoopFilter LPF" (becomes the 2nd sheet / sheetId 5)
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $loopFilter)
$newSheet.Name = "Sheet1"

# ---------------------------------------------------------------------
# 3. "Branch LPF" sheet: add a second-order (Costas) loop filter
#    section below the existing first-order section.
# ---------------------------------------------------------------------
$branch = $wb.Worksheets.Item("Branch LPF")

# Push the existing "Gain" block (old rows 16-22) down by two rows so
# that it starts at row 18, leaving row 17 free for a new header.
$branch.Rows("16:17").Insert()

# New shared strings must be created in the same order they are first
# referenced so they append to sharedStrings.xml in the expected order.
$branch.Range("A29").Value = "B0=b0^2"
$branch.Range("A30").Value = "B1=2b0b1"
$branch.Range("A31").Value = "B2=b1^2"
$branch.Range("A17").Value = "First Order system"
$branch.Range("A26").Value = "Second Order system"
$branch.Range("A33").Value = "A1=2a1"
$branch.Range("A34").Value = "A2=-a1^2"
$branch.Range("A38").Value = "B0_scaled"
$branch.Range("A39").Value = "B1_scaled"
$branch.Range("A40").Value = "B2_scaled"
$branch.Range("A41").Value = "A1_scaled"
$branch.Range("A42").Value = "A2_scaled"

# Header rows, merged + centered
$branch.Range("A17:B17").Merge()
$branch.Range("A17:B17").HorizontalAlignment = -4108
$branch.Range("A26:B26").Merge()
$branch.Range("A26:B26").HorizontalAlignment = -4108

# Second order system, gain
$branch.Range("A27").Value = "Gain"
$branch.Range("B27").Value = 2

# B0, B1, B2 (numerator coefficients)
$branch.Range("B29").Formula = "=B13*B13"
$branch.Range("D29").Value = "Numerator"
$branch.Range("B30").Formula = "=2*B13*B14"
$branch.Range("D30").Value = "Numerator"
$branch.Range("B31").Formula = "=B14*B14"
$branch.Range("D31").Value = "Numerator"

# A1, A2 (denominator coefficients)
$branch.Range("B33").Formula = "=2*B12"
$branch.Range("D33").Value = "Denominator"
$branch.Range("B34").Formula = "=-B12*B12"
$branch.Range("D34").Value = "Denominator"

# Scale bits / scale factor for the second order system
$branch.Range("A36").Value = "scale bits"
$branch.Range("B36").Value = 14
$branch.Range("A37").Value = "scale factor"
$branch.Range("B37").Formula = "=2^B36"

# Scaled coefficients
$branch.Range("B38").Formula = "=ROUND(B29*B27*B37,0)"
$branch.Range("D38").Value = "Numerator"
$branch.Range("B39").Formula = "=ROUND(B30*B27*B37,0)"
$branch.Range("D39").Value = "Numerator"
$branch.Range("B40").Formula = "=ROUND(B31*B27*B37,0)"
$branch.Range("D40").Value = "Numerator"
$branch.Range("B41").Formula = "=ROUND(B33*B37,0)"
$branch.Range("D41").Value = "Denominator"
$branch.Range("B42").Formula = "=ROUND(B34*B37,0)"
$branch.Range("D42").Value = "Denominator"

# ---------------------------------------------------------------------
# 4. Make "Branch LPF" the active sheet/tab, with a view scrolled down
#    to the new content and D43 selected.
# ---------------------------------------------------------------------
$branch.Activate()
$branch.Range("D43").Select()
$excel.ActiveWindow.ScrollRow = 12
